$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Absent" (column H) = 1 for every attendance date row (rows 3-18),
# except row 6 which is handled differently below.
for ($r = 3; $r -le 18; $r++) {
    if ($r -ne 6) {
        $ws.Cells.Item($r, 8).Value = 1
    }
}

# Row 3 also has "Invalid" (column G) = 1
$ws.Cells.Item(3, 7).Value = 1

# Row 6 instead has "Total Attendance Count" (column D) and "Real" (column E) = 1
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
